$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "r585"
$ws.Range("B14").Value = "try"
$ws.Range("C14").Value = "again"
$ws.Range("D14").Value = "2025-10-01 14:51:13"
